# "Generate Report for handoff"
#
# The localization-status report is being regenerated: the handoff
# transform failed for both locales (zh-cn and de-de), so each locale
# sheet's row 2 (the localized markdown file) now reports that failure:
#   - Status ("Ready for handoff" -> "Handoff transform failed")
#   - Latest Handoff File is cleared (value + hyperlink removed)
#   - Latest Handoff Datetime reset to the zero/default date
#   - Handoff Reason flips from "Include" to "Ignored"
# Row 3 (.localization-config, not localized) keeps its values but the
# handoff-datetime column is normalised to the same zero/default date.

$wb = $excel.ActiveWorkbook

$zeroDate = "0001-01-01 00:00:00"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remember the hyperlinks that must survive (everything except the
    # one anchored at C2, the "Latest Handoff File" cell). The sandbox's
    # Hyperlinks.Delete() only works collection-wide, so capture what we
    # need before wiping the collection and then re-add it.
    $keep = @()
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -ne '$C$2') {
            $item = @{
                addr = $addr
                address = $hl.Address
                display = $hl.TextToDisplay
            }
            $keep += , $item
        }
    }

    # --- Cell value updates -------------------------------------------------
    $ws.Range("B2").Value = "Handoff transform failed"
    $ws.Range("C2").ClearContents()
    $ws.Range("D2").Value = $zeroDate
    $ws.Range("G2").Value = $zeroDate
    $ws.Range("H2").Value = "Ignored"

    $ws.Range("D3").Value = $zeroDate
    $ws.Range("G3").Value = $zeroDate
    $ws.Range("H3").Value = "Ignored"

    # --- Hyperlinks: drop C2's link, keep the rest --------------------------
    $ws.Hyperlinks.Delete()
    foreach ($item in $keep) {
        $ws.Hyperlinks.Add($ws.Range($item.addr), $item.address, "", "", $item.display)
    }
}
